# feat: add 2022-Q4 data
#
# Before:  Sheet1 "总计" (summary) + Sheet2 "2022-Q2" (fund holdings for Q2)
# After:   Sheet1 "总计" (summary, now with a Q4 row) + Sheet2 "2022-Q4"
#          (new fund holdings for Q4) + Sheet3 "2022-Q2" (the old Q2 fund
#          holdings sheet, unchanged, moved after the new Q4 sheet)

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)
$q2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. Duplicate the existing "2022-Q2" sheet so its data survives intact
#    as a new sheet placed right after it, then rename both sheets so
#    the original slot becomes "2022-Q4" and the copy keeps "2022-Q2".
# ---------------------------------------------------------------------
$q2.Copy($null, $q2) | Out-Null
$q2Copy = $wb.Worksheets.Item(3)

$q2.Name = "2022-Q4"
$q2Copy.Name = "2022-Q2"

$q4 = $q2

# ---------------------------------------------------------------------
# 2. Replace the (now renamed) "2022-Q4" sheet's contents with the new
#    fund-holdings data for the quarter.
# ---------------------------------------------------------------------
$q4.Cells.Clear() | Out-Null

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("C2").Value = "银华中证全指医药卫生指数增强"
$q4.Range("H2").Value = 7

# Text-like values that look numeric (fund code, fund size, position %,
# etc.) must stay plain text, not be reinterpreted as numbers. Writing
# them as a `="literal"` formula forces text, then copy/paste-values
# bakes that text in as a literal without tagging the destination cell
# with a "quote prefix" style (which a direct numeric-looking
# .Value assignment would otherwise pick up).
function Set-TextValue($range, $text) {
    $helper = $q4.Range("Z100")
    $helper.Formula = '="' + $text + '"'
    $helper.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null  # xlPasteValues
    $excel.CutCopyMode = $false
    $helper.Clear() | Out-Null
}

Set-TextValue $q4.Range("B2") "005112"
Set-TextValue $q4.Range("D2") "1.47"
Set-TextValue $q4.Range("E2") "90.64"
Set-TextValue $q4.Range("F2") "4.99"
Set-TextValue $q4.Range("G2") "0.0734"

# Match the header/first-data-cell styling used elsewhere in the
# workbook (bold, centered, bordered "table header" look) by copying
# the format already used on the summary sheet's header row.
$summary.Range("B1").Copy() | Out-Null
$q4.Range("B1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$summary.Range("A2").Copy() | Out-Null
$q4.Range("A2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Update the "总计" summary sheet: push the existing "2022-Q2" row
#    down to row 3 and insert the new "2022-Q4" row in its old spot.
# ---------------------------------------------------------------------
$summary.Range("A2:D2").Copy() | Out-Null
$summary.Range("A3").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = $false

$summary.Range("A2:D2").Copy() | Out-Null
$summary.Range("A3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$summary.Range("A3").Value = 1

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.07000000000000001
